$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in missing "V" marker in column D for rows 12 and 14 (to match
# the rest of the table, which all have "V" in column D already).
$ws.Range("D12").Value = "V"
$ws.Range("D14").Value = "V"

# Update the active cell / selection on the sheet.
$ws.Range("G11").Select()
